$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Formula = "=PRODUCT(20,30)"
$ws.Range("A2").Select()
